$d = $word.ActiveDocument

# 1) Insert a new empty paragraph at the very start of the document, carrying
#    the same paragraph formatting (spacing/indent/rPr) as the blank paragraph
#    that currently sits right after the picture paragraph. Using raw XML
#    insertion avoids leaving a stray empty run behind (unlike
#    Range.InsertParagraphBefore/After, which clones the neighbouring run's
#    rPr into a new empty <w:r>).
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/><w:ind w:left="720"/><w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr></w:p>'
$d.Range(0, 0).InsertXML($newParaXml)

# 2) Remove the now-redundant blank paragraph that used to live right after
#    the picture (it has effectively been moved above the picture).
$d.Paragraphs.Item(3).Range.Delete()

# 3) Resize/reposition the picture (anchored shape) to its new size/location.
$shape = $d.Shapes.Item(1)
$shape.Left = 403.35614173228345
$shape.Top = 0.8479527559055118
$shape.Width = 64.44338582677166
$shape.Height = 63.43370078740158
